# Generate Report for Handoff
# Inserts a new row (for 397c7b2a-c08e-4bd6-a118-d7cd314ebd53.md) above the
# existing 5c4c0b97-316d-4f0b-a362-7ac1df78037e.md row on each of the three
# sheets (Overview, zh-cn, de-de), then fixes up hyperlinks and resizes the
# tables to include the new row.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c51c107240b0e6f018b5ea1197edcc2d6027889e/e2e/"
$newFile = "397c7b2a-c08e-4bd6-a118-d7cd314ebd53.md"
$oldFile = "5c4c0b97-316d-4f0b-a362-7ac1df78037e.md"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("B2").Value = "e2e\" + $newFile
$ws.Range("C2").Value = ".md"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-09-03 22:43:20"
$ws.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("B2"), $baseUrl + $newFile, "", "", "e2e\" + $newFile)
$null = $ws.Hyperlinks.Add($ws.Range("B3"), $baseUrl + $oldFile, "", "", "e2e\" + $oldFile)

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "397c7b2a-c08e-4bd6-a118-d7cd314ebd53.3004211bb8549588535feb913b5c0caaeee76143.zh-cn.xlf"
$ws.Range("H2").Value = "2016-09-03 22:43:15"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""

$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + $newFile, "", "", $newFile)
$null = $ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + $oldFile, "", "", $oldFile)

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = $newFile
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "e2e"
$ws.Range("E2").Value = "ht"
$ws.Range("F2").Value = "False"
$ws.Range("G2").Value = "397c7b2a-c08e-4bd6-a118-d7cd314ebd53.3004211bb8549588535feb913b5c0caaeee76143.de-de.xlf"
$ws.Range("H2").Value = "2016-09-03 22:43:20"
$ws.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "0001-01-01 00:00:00"
$ws.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "True"
$ws.Range("N2").Value = ""
$ws.Range("O2").Value = "False"
$ws.Range("P2").Value = ""

$ws.Hyperlinks.Delete()
$null = $ws.Hyperlinks.Add($ws.Range("A2"), $baseUrl + $newFile, "", "", $newFile)
$null = $ws.Hyperlinks.Add($ws.Range("A3"), $baseUrl + $oldFile, "", "", $oldFile)

$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:P3"))
